$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style of A2 (bordered/bold/centered header style) for use on new rows in column A
$ws.Range("A2").Copy()

$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = 0.6090525509482942
$ws.Cells.Item(2, 3).Value = 0.9777443151841506
$ws.Cells.Item(2, 4).Value = 150

$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = 0.6158446575419874
$ws.Cells.Item(3, 3).Value = 0.9781730676537722
$ws.Cells.Item(3, 4).Value = 450

$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = 0.5924273823341631
$ws.Cells.Item(4, 3).Value = 0.9785545357999735
$ws.Cells.Item(4, 4).Value = 32

$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 2).Value = 0.6174638305355407
$ws.Cells.Item(5, 3).Value = 0.9790965634947779
$ws.Cells.Item(5, 4).Value = 350

$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = 0.603532451203402
$ws.Cells.Item(6, 3).Value = 0.9777717296365427
$ws.Cells.Item(6, 4).Value = 200

$ws.Cells.Item(7, 1).Value = 5
$ws.Cells.Item(7, 2).Value = 0.6077908466471653
$ws.Cells.Item(7, 3).Value = 0.9786448284056865
$ws.Cells.Item(7, 4).Value = 256

$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = 0.5920415103264995
$ws.Cells.Item(8, 3).Value = 0.9777780505941711
$ws.Cells.Item(8, 4).Value = 25

$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = 0.6208749301584144
$ws.Cells.Item(9, 3).Value = 0.9798951672870251
$ws.Cells.Item(9, 4).Value = 512

$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = 0.6254428465434397
$ws.Cells.Item(10, 3).Value = 0.9801289712231603
$ws.Cells.Item(10, 4).Value = 1024

$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = 0.6092459848759443
$ws.Cells.Item(11, 3).Value = 0.9780775366505514
$ws.Cells.Item(11, 4).Value = 400

$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = 0.6085730073322152
$ws.Cells.Item(12, 3).Value = 0.9767224112683015
$ws.Cells.Item(12, 4).Value = 300

$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 2).Value = 0.6221821314463949
$ws.Cells.Item(13, 3).Value = 0.9783465401317034
$ws.Cells.Item(13, 4).Value = 500

$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 2).Value = 0.6057427544415619
$ws.Cells.Item(14, 3).Value = 0.9781840931880694
$ws.Cells.Item(14, 4).Value = 250

$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(15, 2).Value = 0.5967852939947144
$ws.Cells.Item(15, 3).Value = 0.9773693057641633
$ws.Cells.Item(15, 4).Value = 16

$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = 0.5952083146592653
$ws.Cells.Item(16, 3).Value = 0.9776092170085141
$ws.Cells.Item(16, 4).Value = 128

# Apply the header-like bordered style (style index 1) to the new column-A cells in rows 11-16
$ws.Range("A11:A16").PasteSpecial(-4122)

$excel.CutCopyMode = 0